$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1426875065163276
$ws.Range("C2").Value = 1.09163566462907
$ws.Range("D2").Value = 0.1621774786541525
$ws.Range("B3").Value = 0.06780569901143271
$ws.Range("C3").Value = 5.467029463578041
$ws.Range("D3").Value = 0.2460062864740273
$ws.Range("B4").Value = 0.09454485513604982
$ws.Range("C4").Value = 8.927750571830922
$ws.Range("D4").Value = 0.1187085183396777
$ws.Range("B5").Value = 0.1341618055095846
$ws.Range("C5").Value = 12.10800921252872
$ws.Range("D5").Value = 0.1240484889818186
$ws.Range("B6").Value = 0.09713028466478685
$ws.Range("C6").Value = 14.25127202409365
$ws.Range("D6").Value = 0.128325096430879
$ws.Range("B7").Value = 0.08319026540728862
$ws.Range("C7").Value = 16.38686896093605
$ws.Range("D7").Value = 0.2352949081865892
$ws.Range("B8").Value = 0.1437298375993895
$ws.Range("C8").Value = 18.56989698139986
$ws.Range("D8").Value = 0.1377585931658061
$ws.Range("B9").Value = 0.1436556778735269
$ws.Range("C9").Value = 20.24744671788237
$ws.Range("D9").Value = 0.1580391792726555
$ws.Range("B10").Value = 0.1278563196084907
$ws.Range("C10").Value = 21.33210122606258
$ws.Range("D10").Value = 0.2539114923373028
$ws.Range("B11").Value = 0.0796381373571928
$ws.Range("C11").Value = 21.23520642445054
$ws.Range("D11").Value = 0.22040003565719
$ws.Range("B12").Value = 0.08469928759339618
$ws.Range("C12").Value = 21.43165087645074
$ws.Range("D12").Value = 0.2067868352540644
$ws.Range("B13").Value = 0.130587601614634
$ws.Range("C13").Value = 21.89252627870857
$ws.Range("D13").Value = 0.2861840985946398
$ws.Range("B14").Value = 0.1246631449503537
$ws.Range("C14").Value = 21.27424046679059
$ws.Range("D14").Value = 0.2119706330389948
$ws.Range("B15").Value = 0.1267055264446616
$ws.Range("C15").Value = 19.44357778612015
$ws.Range("D15").Value = 0.2664298416002698
$ws.Range("B16").Value = 0.09746172721490505
$ws.Range("C16").Value = 18.49109689988184
$ws.Range("D16").Value = 0.1919594135899338
$ws.Range("B17").Value = 0.09386597917387131
$ws.Range("C17").Value = 16.25194462883666
$ws.Range("D17").Value = 0.2617381314950782
$ws.Range("B18").Value = 0.1207858929496632
$ws.Range("C18").Value = 13.73795313474392
$ws.Range("D18").Value = 0.1457683181667988
$ws.Range("B19").Value = 0.1333385376525511
$ws.Range("C19").Value = 11.35863696581246
$ws.Range("D19").Value = 0.2860747630596674
$ws.Range("B20").Value = 0.08966144511699978
$ws.Range("C20").Value = 7.518375554057935
$ws.Range("D20").Value = 0.1147328330642163
$ws.Range("B21").Value = 0.07301333228164658
$ws.Range("C21").Value = 3.783536648721564
$ws.Range("D21").Value = 0.2323711876133207
$ws.Range("B22").Value = 0.09656121264505674
$ws.Range("C22").Value = -0.2987977553288744
$ws.Range("D22").Value = 0.2553141745934079
$ws.Range("B23").Value = 0.1221450509918528
$ws.Range("C23").Value = -4.702016228327486
$ws.Range("D23").Value = 0.1735526591591772
$ws.Range("B24").Value = 0.1455195460935777
$ws.Range("C24").Value = -10.08396686336512
$ws.Range("D24").Value = 0.2612427949436753
$ws.Range("B25").Value = 0.1438335567502133
$ws.Range("C25").Value = -15.26967918809291
$ws.Range("D25").Value = 0.2546023436770116
$ws.Range("B26").Value = 0.06139636678659802
$ws.Range("C26").Value = -20.87635839077612
$ws.Range("D26").Value = 0.1417412915639149
